$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The practice-trial table (rows 2-10) holds 3 blocks of 3 rows each
# (rot / blau / gelb), and within every block the "congruent" trial
# (word == color, column D = 1) gets duplicated immediately below
# itself so it is shown twice (speed vs. accuracy practice).
#
# Congruent rows in the original sheet: row 2 (rot/red/r/1),
# row 6 (blau/blue/b/1), row 10 (gelb/yellow/g/1).
#
# Insert from the bottom up so earlier row numbers stay valid.

# --- gelb block: duplicate row 10 (gelb/yellow/g/1) -> new row 11 ---
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "gelb"
$ws.Range("B11").Value = "yellow"
$ws.Range("C11").Value = "g"
$ws.Range("D11").Value = 1

# --- blau block: duplicate row 6 (blau/blue/b/1) -> new row 7 ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "blau"
$ws.Range("B7").Value = "blue"
$ws.Range("C7").Value = "b"
$ws.Range("D7").Value = 1

# --- rot block: duplicate row 2 (rot/red/r/1) -> new row 3 ---
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "rot"
$ws.Range("B3").Value = "red"
$ws.Range("C3").Value = "r"
$ws.Range("D3").Value = 1
# Row-insert copies the formatting of the row above (row 2, which is
# fully styled s="2"); only column A keeps that styling in the final
# sheet, so clear it back off of B3:C3.
$ws.Range("B3:C3").ClearFormats()

$ws.Range("F5").Select()
